$wb = $excel.ActiveWorkbook

# Update the "Time" column header to "Timestamp" on both manager sheets.
$ws1 = $wb.Worksheets.Item("Manager1")
$ws1.Range("D1").Value = "Timestamp"
$ws1.Range("F3").Select()

$ws2 = $wb.Worksheets.Item("Manager2")
$ws2.Range("D1").Value = "Timestamp"
$ws2.Range("D1").Select()
